$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the H1 title.
#    Borrow the run/paragraph shape (leading empty run + bold run) from
#    the existing bold "title" paragraph near the end of the doc so the
#    new paragraph doesn't inherit the Heading1 style of paragraph 1.
# ---------------------------------------------------------------------
$titlePara = $d.Paragraphs.Item(1)
$titlePara.Range.InsertParagraphAfter() | Out-Null

$count = $d.Paragraphs.Count
$boldSource = $d.Paragraphs.Item($count - 1)
$boldFormatted = $boldSource.Range.FormattedText

$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Style = "Normal"
$metaPara.Range.FormattedText = $boldFormatted

$metaRange = $d.Paragraphs.Item(2).Range
$metaRange.MoveEnd(1, -1) | Out-Null
$metaRange.Text = "Meta description"

$metaRest = $d.Range($metaRange.End, $metaRange.End)
$metaRest.InsertAfter(": Read our review of Book of Oz Lock 'n Spin, a unique online slot game with a touch of magic and customizable Lock 'n Spin feature. Play for free now.")
$metaRest.Font.Bold = 0

# ---------------------------------------------------------------------
# 2) Remove the duplicated bold title paragraph near the end of the doc
#    and replace the italic meta-description paragraph's text with the
#    DALLE image prompt.
# ---------------------------------------------------------------------
$n = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs.Item($n - 1)
$dupTitlePara.Range.Delete() | Out-Null

$n2 = $d.Paragraphs.Count
$italicPara = $d.Paragraphs.Item($n2)
$italicRange = $italicPara.Range
$italicRange.MoveEnd(1, -1) | Out-Null

$quote = [char]34
$lsquo = [char]0x2018
$dallePrompt = "Prompt for DALLE: Create a cartoon-style feature image for " + $quote + "Book of Oz Lock " + $lsquo + "N Spins" + $quote + " that features a happy Maya warrior with glasses. The background should be green and the warrior should be holding a magic book with the game's title on it. The warrior should be surrounded by symbols from the game, such as playing card symbols and magic filters in the shape of flowers, hearts, spades, and diamonds. Use bright colors and make the image dynamic and engaging to attract potential players."
$italicRange.Text = $dallePrompt

Write-Output "done"
